$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "SP.RUR.TOTL:IMN"
$ws.Range("C1").Value = "SP.URB.TOTL:IMN"
$ws.Range("A2").Value = "SP.POP.TOTL:IMN:cor-value"
$ws.Range("A3").Value = "SP.POP.TOTL:IMN:p-value"
$ws.Range("A4").Value = "SP.RUR.TOTL:IMN:cor-value"
$ws.Range("A5").Value = "SP.RUR.TOTL:IMN:p-value"

$ws.Range("B2").Value = 0.9995941726775228
$ws.Range("C2").Value = 0.9997299486206004
$ws.Range("B3").Value = [double]"6.444095796974465e-20"
$ws.Range("C3").Value = [double]"5.596532955662555e-21"
$ws.Range("C4").Value = 0.998662242387819
$ws.Range("C5").Value = [double]"8.251089475053815e-17"

$headerRange = $ws.Range("B1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$labelRange = $ws.Range("A2:A5")
$labelRange.Font.Bold = $true
$labelRange.HorizontalAlignment = -4108
$labelRange.VerticalAlignment = -4160
$labelRange.Borders.LineStyle = 1
$labelRange.Borders.Weight = 2
